# LMS-2523 Update BaSynthec Validation
# Rename the two hybridization-series header labels on the "openbis-data"
# sheet (D1/E1) from "MGP253"/"MGP776" to "JJS-MGP253"/"JJS-MGP776", then
# leave the workbook focused on that sheet (as the author last left it),
# with the selection parked at I16.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("openbis-data")

$wsData.Range("D1").Value = "JJS-MGP253"
$wsData.Range("E1").Value = "JJS-MGP776"

$wsData.Activate() | Out-Null
$wsData.Range("I16").Select() | Out-Null
